$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.828.69'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.07%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.740.54'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.94%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.08%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5163'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.77%  '

$ws.Range("E8").Value = '  +6.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.14'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.08%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06089'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.74%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.733.80'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.37%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06965'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.37%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.18'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.94%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6328'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.09%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.494'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.97%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '76.41'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.01%  '

$ws.Range("E18").Value = '  -0.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '25.848.36'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.01%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.42'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.90%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006563'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.67%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.962.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.80%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.068'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.421'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.34%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.103'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.33%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '137.44'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.31%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.505'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.14%  '

$ws.Range("E28").Value = '  +0.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '14.95'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.37%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '102.53'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.10%  '

$ws.Range("E31").Value = '  +0.44%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.605'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.31%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.404'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04392'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.622'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.26%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9673'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.56%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5996'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.667'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.31%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01547'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.04%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9997'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.18%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.898'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.98%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '100.68'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.60%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.3817'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.06%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7226'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.77%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.881'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.05455'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.242'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.88%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1097'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.88%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.69'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.58%  '

$ws.Range("E50").Value = '  +0.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.486'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.47%  '
